$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1): columns E..K are new, give them proper column
#     names (matching the header convention used on every other sheet) ---
$ws.Cells.Item(1, 5).Value  = "property_category"
$ws.Cells.Item(1, 6).Value  = "category"
$ws.Cells.Item(1, 7).Value  = "date"
$ws.Cells.Item(1, 8).Value  = "legislator_name"
$ws.Cells.Item(1, 9).Value  = "legislator_id"
$ws.Cells.Item(1, 10).Value = "source_file"
$ws.Cells.Item(1, 11).Value = "index"

# Also fix the two pre-existing header cells that were wrongly populated
# with data values instead of field names (C1/D1 — B1 already happens to
# hold the right text, "company" is still the bank/company name there).
$ws.Cells.Item(1, 2).Value = "company"
$ws.Cells.Item(1, 3).Value = "name"
$ws.Cells.Item(1, 4).Value = "owner"

# Match the bold / centered / bordered look already used by B1:D1 (and by
# every other sheet's header row) for the newly added header cells.
foreach ($col in 5..11) {
    $c = $ws.Cells.Item(1, $col)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
    $c.Borders.LineStyle = 1
}

# --- Row 2 (record #91) ---
$ws.Cells.Item(2, 2).Value  = "中華郵政"
$ws.Cells.Item(2, 3).Value  = "六六金順"
$ws.Cells.Item(2, 4).Value  = "楊瓊瓔"
$ws.Cells.Item(2, 5).Value  = "insurance"
$ws.Cells.Item(2, 6).Value  = "normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value  = "2011-12-27"
$ws.Cells.Item(2, 8).Value  = "楊瓊瓔"
$ws.Cells.Item(2, 9).Value  = 854
$ws.Cells.Item(2, 10).Value = "tmpd1401"
$ws.Cells.Item(2, 11).Value = 91

# --- Row 3 (record #92) ---
$ws.Cells.Item(3, 2).Value  = "中華郵政"
$ws.Cells.Item(3, 3).Value  = "吉利保險"
$ws.Cells.Item(3, 4).Value  = "楊瓊壤"
$ws.Cells.Item(3, 5).Value  = "insurance"
$ws.Cells.Item(3, 6).Value  = "normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value  = "2011-12-27"
$ws.Cells.Item(3, 8).Value  = "楊瓊瓔"
$ws.Cells.Item(3, 9).Value  = 854
$ws.Cells.Item(3, 10).Value = "tmpd1401"
$ws.Cells.Item(3, 11).Value = 92

# The date column was forced to Text via NumberFormat "@" so it stores the
# literal string "2011-12-27" (like every other "date" column in this
# workbook) instead of being auto-converted to a date serial number.
# Drop back to the plain/default look used by the rest of the data rows.
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 7).Style = "Normal"
